$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 247
$ws.Cells.Item(4, 6).Value = 127
$ws.Cells.Item(5, 6).Value = 357
$ws.Cells.Item(6, 6).Value = 518
$ws.Cells.Item(9, 6).Value = 260
$ws.Cells.Item(12, 6).Value = 587
$ws.Cells.Item(13, 6).Value = 734
$ws.Cells.Item(14, 6).Value = 1489
$ws.Cells.Item(15, 6).Value = 1489
$ws.Cells.Item(18, 6).Value = 1339
$ws.Cells.Item(20, 6).Value = 253
$ws.Cells.Item(24, 6).Value = 6440
$ws.Cells.Item(25, 6).Value = 4785
$ws.Cells.Item(27, 6).Value = 487
$ws.Cells.Item(28, 6).Value = 192
$ws.Cells.Item(29, 6).Value = 117
$ws.Cells.Item(32, 6).Value = 1249
$ws.Cells.Item(33, 6).Value = 185
$ws.Cells.Item(34, 6).Value = 36
$ws.Cells.Item(37, 6).Value = 1332
$ws.Cells.Item(38, 6).Value = 225
$ws.Cells.Item(40, 6).Value = 139
$ws.Cells.Item(41, 6).Value = 57

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(15, 6).Value = 235

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 2429
$ws.Cells.Item(4, 6).Value = 173
$ws.Cells.Item(5, 6).Value = 40

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 247
$ws.Cells.Item(7, 6).Value = 173
$ws.Cells.Item(8, 6).Value = 40
$ws.Cells.Item(9, 6).Value = 357
$ws.Cells.Item(10, 6).Value = 518
$ws.Cells.Item(13, 6).Value = 260
$ws.Cells.Item(17, 6).Value = 587
$ws.Cells.Item(18, 6).Value = 734
$ws.Cells.Item(19, 6).Value = 1489
$ws.Cells.Item(20, 6).Value = 1489
$ws.Cells.Item(23, 6).Value = 1339
$ws.Cells.Item(25, 6).Value = 253
$ws.Cells.Item(31, 6).Value = 6440
$ws.Cells.Item(32, 6).Value = 4785
$ws.Cells.Item(34, 6).Value = 192
$ws.Cells.Item(36, 6).Value = 1249
$ws.Cells.Item(37, 6).Value = 185
$ws.Cells.Item(38, 6).Value = 36
$ws.Cells.Item(44, 6).Value = 1332
$ws.Cells.Item(45, 6).Value = 225
$ws.Cells.Item(46, 6).Value = 139
$ws.Cells.Item(47, 6).Value = 57
$ws.Cells.Item(49, 6).Value = 235

